$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, $cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

function Clear-CellValue($ws, $cellRef) {
    $ws.Range($cellRef).Value = $null
}

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 31
Set-CellValue $ws "H31" 1915.3334
Set-CellValue $ws "I31" 108.5
Set-CellValue $ws "J31" 10949.5
Set-CellValue $ws "K31" 325.5
Set-CellValue $ws "L31" 32848.5
Set-CellValue $ws "M31" -95.5
Set-CellValue $ws "N31" -33308.5

# Row 76
Set-CellValue $ws "H76" 4866.091
Set-CellValue $ws "I76" 4947.6665
Set-CellValue $ws "K76" 4947.6665
Set-CellValue $ws "M76" -4632.6665

# Row 79
Set-CellValue $ws "H79" 4866.091
Set-CellValue $ws "I79" 4947.6665
Set-CellValue $ws "K79" 4947.6665
Set-CellValue $ws "M79" -3855.6665

# Row 106
Set-CellValue $ws "H106" 10405.5
Set-CellValue $ws "I106" 10405.5
Set-CellValue $ws "K106" 10405.5
Set-CellValue $ws "M106" -9774.5

# Row 116
Set-CellValue $ws "H116" 6169.3335
Set-CellValue $ws "I116" 5775.4
Set-CellValue $ws "J116" 6661.75
Set-CellValue $ws "K116" 5775.4
Set-CellValue $ws "L116" 6661.75
Set-CellValue $ws "M116" -2333.4
Set-CellValue $ws "N116" -13545.75

# Row 138
Set-CellValue $ws "H138" 10925.867
Set-CellValue $ws "J138" 11013
Set-CellValue $ws "L138" 33039
Set-CellValue $ws "N138" -43319

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
Set-CellValue $ws "H32" 7424.1665
Set-CellValue $ws "I32" 7424.1665
Set-CellValue $ws "K32" 7424.1665
Set-CellValue $ws "M32" -7137.1665

# Row 41
Set-CellValue $ws "H41" 5363.4
Set-CellValue $ws "I41" 5892.778
Set-CellValue $ws "K41" 5892.778
Set-CellValue $ws "M41" -5478.778

# Row 44
Set-CellValue $ws "H44" 49999
Set-CellValue $ws "J44" 49999
Set-CellValue $ws "L44" 49999
Set-CellValue $ws "N44" -50975

# Row 45
Set-CellValue $ws "H45" 6998.392
Set-CellValue $ws "I45" 14082.5625
Set-CellValue $ws "K45" 14082.5625
Set-CellValue $ws "M45" -13705.5625

# Row 92
Set-CellValue $ws "H92" 45798.4
Set-CellValue $ws "J92" 45798.4
Set-CellValue $ws "L92" 45798.4
Set-CellValue $ws "N92" -50790.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
Set-CellValue $ws "H20" 1838.36
Set-CellValue $ws "I20" 1620.909
Set-CellValue $ws "K20" 1620.909
Set-CellValue $ws "M20" -1373.909

# Row 92
Set-CellValue $ws "H92" 88333.25
Set-CellValue $ws "J92" 88333.25
Set-CellValue $ws "L92" 88333.25
Set-CellValue $ws "N92" -93325.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
Set-CellValue $ws "H16" 61843
Set-CellValue $ws "I16" 862.5
Set-CellValue $ws "J16" 183804
Set-CellValue $ws "K16" 862.5
Set-CellValue $ws "L16" 183804
Set-CellValue $ws "M16" -575.5
Set-CellValue $ws "N16" -184378

# Row 31
Set-CellValue $ws "H31" 7182.9165
Set-CellValue $ws "J31" 14000
Set-CellValue $ws "L31" 14000
Set-CellValue $ws "N31" -14590

# Row 34
Set-CellValue $ws "H34" 7182.9165
Set-CellValue $ws "J34" 14000
Set-CellValue $ws "L34" 14000
Set-CellValue $ws "N34" -14404

# Row 58
Set-CellValue $ws "H58" 4915.125
Set-CellValue $ws "I58" 5151.875
Set-CellValue $ws "J58" 4441.625
Set-CellValue $ws "K58" 5151.875
Set-CellValue $ws "L58" 4441.625
Set-CellValue $ws "M58" -4948.875
Set-CellValue $ws "N58" -4847.625

# Row 62
Set-CellValue $ws "H62" 1597.4584
Set-CellValue $ws "I62" 1492.3334
Set-CellValue $ws "K62" 1492.3334
Set-CellValue $ws "M62" -868.3334

# Row 65
Set-CellValue $ws "H65" 1597.4584
Set-CellValue $ws "I65" 1492.3334
Set-CellValue $ws "K65" 7461.666999999999
Set-CellValue $ws "M65" -4341.666999999999

# Row 107
Set-CellValue $ws "H107" 1109.54
Set-CellValue $ws "I107" 860.725
Set-CellValue $ws "K107" 860.725
Set-CellValue $ws "M107" 1059.275

# Row 113
Set-CellValue $ws "H113" 61843
Set-CellValue $ws "I113" 862.5
Set-CellValue $ws "J113" 183804
Set-CellValue $ws "K113" 862.5
Set-CellValue $ws "L113" 183804
Set-CellValue $ws "M113" 1307.5
Set-CellValue $ws "N113" -188144

# Row 134
Set-CellValue $ws "H134" 3559.8
Set-CellValue $ws "J134" 2600
Set-CellValue $ws "L134" 7800
Set-CellValue $ws "N134" -12870

# Row 136
Set-CellValue $ws "H136" 4915.125
Set-CellValue $ws "I136" 5151.875
Set-CellValue $ws "J136" 4441.625
Set-CellValue $ws "K136" 15455.625
Set-CellValue $ws "L136" 13324.875
Set-CellValue $ws "M136" -12905.625
Set-CellValue $ws "N136" -18424.875

# Row 137
Set-CellValue $ws "H137" 89807.734
Set-CellValue $ws "J137" 89807.734
Set-CellValue $ws "L137" 89807.734
Set-CellValue $ws "N137" -100007.734

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 122
Set-CellValue $ws "H122" 2498
Set-CellValue $ws "J122" 2498
Set-CellValue $ws "L122" 22482
Set-CellValue $ws "N122" -27382

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
Set-CellValue $ws "H132" 2999
Set-CellValue $ws "I132" 2998.5
Set-CellValue $ws "K132" 8995.5
Set-CellValue $ws "M132" -6465.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
Set-CellValue $ws "H46" 2437.6667
Set-CellValue $ws "I46" 2200
Set-CellValue $ws "K46" 2200
Set-CellValue $ws "M46" -2012

# Row 94
Set-CellValue $ws "H94" 0
Set-CellValue $ws "J94" 0
Set-CellValue $ws "L94" 0
Clear-CellValue $ws "N94"

# Row 100
Set-CellValue $ws "H100" 4371.8
Set-CellValue $ws "I100" 4209.3076
Set-CellValue $ws "K100" 4209.3076
Set-CellValue $ws "M100" -3668.3076

# Row 115
Set-CellValue $ws "H115" 125249
Set-CellValue $ws "J115" 125249
Set-CellValue $ws "L115" 125249
Set-CellValue $ws "N115" -127599

# Row 122
Set-CellValue $ws "H122" 5233.878
Set-CellValue $ws "I122" 5078.1816
Set-CellValue $ws "J122" 5876.125
Set-CellValue $ws "K122" 15234.5448
Set-CellValue $ws "L122" 17628.375
Set-CellValue $ws "M122" -12784.5448
Set-CellValue $ws "N122" -22528.375

# Row 132
Set-CellValue $ws "H132" 18199.928
Set-CellValue $ws "I132" 13679.9
Set-CellValue $ws "J132" 29500
Set-CellValue $ws "K132" 41039.7
Set-CellValue $ws "L132" 88500
Set-CellValue $ws "M132" -38509.7
Set-CellValue $ws "N132" -93560

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
Set-CellValue $ws "H96" 1798.5555
Set-CellValue $ws "I96" 1916.5
Set-CellValue $ws "J96" 1562.6666
Set-CellValue $ws "K96" 1916.5
Set-CellValue $ws "L96" 1562.6666
Set-CellValue $ws "M96" -543.5
Set-CellValue $ws "N96" -4308.6666

# Row 110
Set-CellValue $ws "H110" 0
Set-CellValue $ws "J110" 0
Set-CellValue $ws "L110" 0
Clear-CellValue $ws "N110"

# Row 126
Set-CellValue $ws "H126" 1807.9286
Set-CellValue $ws "I126" 1652.6666
Set-CellValue $ws "K126" 4957.9998
Set-CellValue $ws "M126" -2487.9998

# Row 132
Set-CellValue $ws "H132" 5737.3706
Set-CellValue $ws "I132" 5141.7617
Set-CellValue $ws "K132" 15425.2851
Set-CellValue $ws "M132" -12895.2851

# Row 136
Set-CellValue $ws "H136" 15181.809
Set-CellValue $ws "I136" 18477.229
Set-CellValue $ws "K136" 55431.687
Set-CellValue $ws "M136" -52881.687
